$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 36
$data = @(
    @('35', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'admin_dashboard', '2025-07-03 16:16:11'),
    @('36', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:16:12'),
    @('37', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:16:12'),
    @('38', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:16:12'),
    @('39', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'home', '2025-07-03 16:16:15'),
    @('40', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:16:16'),
    @('41', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:16:16'),
    @('42', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:16:16'),
    @('43', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'unknown', '2025-07-03 16:16:38'),
    @('44', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'home', '2025-07-03 16:18:03'),
    @('45', '127.0.0.1', 'Mozilla/5.0 (X11; Linux x86_64) AppleWebKit/537.36 (KHTML, like Gecko) HeadlessChrome/118.0.5993.88 ...', 'home', '2025-07-03 16:18:03'),
    @('46', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'home', '2025-07-03 16:18:03'),
    @('47', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'home', '2025-07-03 16:18:04'),
    @('48', '127.0.0.1', 'Mozilla/5.0 (X11; Linux x86_64) AppleWebKit/537.36 (KHTML, like Gecko) HeadlessChrome/118.0.5993.88 ...', 'static', '2025-07-03 16:18:04'),
    @('49', '127.0.0.1', 'Mozilla/5.0 (X11; Linux x86_64) AppleWebKit/537.36 (KHTML, like Gecko) HeadlessChrome/118.0.5993.88 ...', 'static', '2025-07-03 16:18:04'),
    @('50', '127.0.0.1', 'Mozilla/5.0 (X11; Linux x86_64) AppleWebKit/537.36 (KHTML, like Gecko) HeadlessChrome/118.0.5993.88 ...', 'static', '2025-07-03 16:18:04'),
    @('51', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:04'),
    @('52', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:04'),
    @('53', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:04'),
    @('54', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:04'),
    @('55', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:04'),
    @('56', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:04'),
    @('57', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'home', '2025-07-03 16:18:35'),
    @('58', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'home', '2025-07-03 16:18:36'),
    @('59', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:37'),
    @('60', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:37'),
    @('61', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:37'),
    @('62', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'vote_comment', '2025-07-03 16:18:42'),
    @('63', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'vote_comment', '2025-07-03 16:18:44'),
    @('64', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'vote_comment', '2025-07-03 16:18:45'),
    @('65', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'admin_dashboard', '2025-07-03 16:18:51'),
    @('66', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:51'),
    @('67', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:52'),
    @('68', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:18:52'),
    @('69', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'delete_topic', '2025-07-03 16:19:00'),
    @('70', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'admin_dashboard', '2025-07-03 16:19:01'),
    @('71', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:19:01'),
    @('72', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:19:02'),
    @('73', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:19:02'),
    @('74', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'admin_dashboard', '2025-07-03 16:19:07'),
    @('75', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:19:08'),
    @('76', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:19:08'),
    @('77', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'static', '2025-07-03 16:19:08'),
    @('78', '31.94.64.109', 'Mozilla/5.0 (iPhone; CPU iPhone OS 18_5 like Mac OS X) AppleWebKit/605.1.15 (KHTML, like Gecko) Mobi...', 'delete_topic', '2025-07-03 16:19:12')
)

$endRow = $startRow + $data.Count - 1
$rng = $ws.Range("A" + $startRow + ":E" + $endRow)
$arr = New-Object "object[,]" $data.Count,5
for ($i = 0; $i -lt $data.Count; $i++) {
    $arr[$i,0] = [int]$data[$i][0]
    $arr[$i,1] = $data[$i][1]
    $arr[$i,2] = $data[$i][2]
    $arr[$i,3] = $data[$i][3]
    $arr[$i,4] = $data[$i][4]
}
$rng.Value = $arr

